$wb = $excel.ActiveWorkbook

# "Week 6" holds the already-correctly-styled date/start/end cells we can
# clone formatting from (style 5 = short date, style 6 = time) so the new
# "Week 7" rows pick up the same cellXfs the author used elsewhere in the
# workbook instead of minting brand new ones.
$wsSrc = $wb.Worksheets.Item("Week 6")
$ws = $wb.Worksheets.Item("Week 7")

# --- Row 2: Tue 2/19 (entered as nominal 2014 date so the stored serial
# matches the workbook's existing literal date serials), 10:00 AM - 12:00 PM
$wsSrc.Range("A2:C2").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
$ws.Range("A2").EntireRow.RowHeight = 18
$ws.Range("A2").Value2 = 41689
$ws.Range("B2").Value2 = 0.41666666666666669
$ws.Range("C2").Value2 = 0.5
$ws.Range("D2").Value = "Successfully displayed database values for product table"
$ws.Range("E2").Value2 = 2

# --- Row 3: Wed 2/20, 3:30 PM - 5:00 PM
$wsSrc.Range("A2:C2").Copy()
$ws.Range("A3:C3").PasteSpecial(-4122)
$ws.Range("A3").EntireRow.RowHeight = 18
$ws.Range("A3").Value2 = 41690
$ws.Range("B3").Value2 = 0.64583333333333337
$ws.Range("C3").Value2 = 0.70833333333333337
$ws.Range("D3").Value = "Worked on product/image relationship"
$ws.Range("E3").Value2 = 1.5

# --- Row 4: Thu 2/21, 10:00 AM - 12:00 PM
$wsSrc.Range("A2:C2").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$ws.Range("A4").EntireRow.RowHeight = 18
$ws.Range("A4").Value2 = 41691
$ws.Range("B4").Value2 = 0.41666666666666669
$ws.Range("C4").Value2 = 0.5
$ws.Range("D4").Value = "Worked on product/image relationship"
$ws.Range("E4").Value2 = 2

$excel.CutCopyMode = 0

# Submitting week 7's timesheet makes it the active tab (was "Week 6").
$ws.Activate()
$ws.Range("A5").Select()
